$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Unprotect()

# Update the confidential disclosure date string (shared string used by A80)
$ws.Range("A80").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-03-31 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) columns for holdings rows 2-77
$ws.Range("D2").Value = 0.07373371597954634
$ws.Range("E2").Value = 0.01876563803169318
$ws.Range("D3").Value = 0.04509315464790287
$ws.Range("E3").Value = 0.01269601248981278
$ws.Range("D4").Value = 0.03535946779286687
$ws.Range("E4").Value = 0.01690748328660785
$ws.Range("D5").Value = 0.03335937156391044
$ws.Range("E5").Value = -0.0008022652194431545
$ws.Range("D6").Value = 0.03191967617480256
$ws.Range("E6").Value = -0.01456499223200414
$ws.Range("D7").Value = 0.03020393749604581
$ws.Range("E7").Value = 0.007842773165499528
$ws.Range("D8").Value = 0.03084826809033046
$ws.Range("E8").Value = -0.003999757590448905
$ws.Range("D9").Value = 0.02921488486815403
$ws.Range("E9").Value = -0.01053487741585646
$ws.Range("D10").Value = 0.0263014007019991
$ws.Range("E10").Value = 0.006359704909692176
$ws.Range("D11").Value = 0.02798757209915204
$ws.Range("E11").Value = -0.000147655961609261
$ws.Range("D12").Value = 0.02373145784605945
$ws.Range("E12").Value = -0.005443863526114257
$ws.Range("D13").Value = 0.02416911264996541
$ws.Range("E13").Value = -0.007694280584765512
$ws.Range("D14").Value = 0.01970804449477388
$ws.Range("E14").Value = 0.02273612317646712
$ws.Range("D15").Value = 0.0186357508677579
$ws.Range("E15").Value = 0.02872228088701156
$ws.Range("D16").Value = 0.02113942562333094
$ws.Range("E16").Value = -0.001158972377824918
$ws.Range("D17").Value = 0.01869832920752503
$ws.Range("E17").Value = 0.0006630322675702516
$ws.Range("D18").Value = 0.01823253382944707
$ws.Range("E18").Value = -0.007646976287357998
$ws.Range("D19").Value = 0.01533931874818525
$ws.Range("E19").Value = -0.01517022402540125
$ws.Range("D20").Value = 0.0141850632390529
$ws.Range("E20").Value = 0.01339076692574204
$ws.Range("D21").Value = 0.01558554877219035
$ws.Range("E21").Value = 0.02267361111111099
$ws.Range("D22").Value = 0.01409783729848127
$ws.Range("E22").Value = -0.006002233389168055
$ws.Range("D23").Value = 0.01298913803355045
$ws.Range("E23").Value = 0.02645214071448043
$ws.Range("D24").Value = 0.01542054267818174
$ws.Range("E24").Value = -0.0149690855841198
$ws.Range("D25").Value = 0.01470236767509956
$ws.Range("E25").Value = -0.003988649748367101
$ws.Range("D26").Value = 0.01192097693840589
$ws.Range("E26").Value = 0.01929743471227141
$ws.Range("D27").Value = 0.01231691284913376
$ws.Range("E27").Value = 0.003275283591627876
$ws.Range("D28").Value = 0.01223539373828932
$ws.Range("E28").Value = 0.003606711619883907
$ws.Range("D29").Value = 0.01190129821520868
$ws.Range("E29").Value = -0.01582393597671872
$ws.Range("D30").Value = 0.01126724975379458
$ws.Range("E30").Value = 0.03961156909319552
$ws.Range("D31").Value = 0.01286988497097536
$ws.Range("E31").Value = -0.002064220183486287
$ws.Range("D32").Value = 0.01349655391119052
$ws.Range("E32").Value = 0.001410668591299791
$ws.Range("D33").Value = 0.01110544145230552
$ws.Range("E33").Value = 0.0108755842027155
$ws.Range("D34").Value = 0.01181800801927649
$ws.Range("E34").Value = -0.007509813961426803
$ws.Range("D35").Value = 0.009381142528957982
$ws.Range("E35").Value = 0.05083225826751825
$ws.Range("D36").Value = 0.01124368448276592
$ws.Range("E36").Value = -0.008278457196613243
$ws.Range("D37").Value = 0.01105747456451232
$ws.Range("E37").Value = -0.01441537640149482
$ws.Range("D38").Value = 0.01030623930645882
$ws.Range("E38").Value = -0.009881139911213088
$ws.Range("D39").Value = 0.008976646373639329
$ws.Range("E39").Value = 0.01644159943879342
$ws.Range("D40").Value = 0.009426108411463605
$ws.Range("E40").Value = 0.0007515657620043026
$ws.Range("D41").Value = 0.008961395363161491
$ws.Range("E41").Value = 0.0190607947121666
$ws.Range("D42").Value = 0.009157001871741757
$ws.Range("E42").Value = -0.0004512974802557412
$ws.Range("D43").Value = 0.009844084492172346
$ws.Range("E43").Value = 0.001689189189189255
$ws.Range("D44").Value = 0.009718337450942173
$ws.Range("E44").Value = -0.005386250885896504
$ws.Range("D45").Value = 0.009380552167262066
$ws.Range("E45").Value = -0.01076182384593616
$ws.Range("D46").Value = 0.009593082377791932
$ws.Range("E46").Value = 0.003323179174743673
$ws.Range("D47").Value = 0.009001835139331756
$ws.Range("E47").Value = -0.01757607555089169
$ws.Range("D48").Value = 0.007300806306164922
$ws.Range("E48").Value = 0.007115902964959453
$ws.Range("D49").Value = 0.008366507560909831
$ws.Range("E49").Value = -0.008655666756829827
$ws.Range("D50").Value = 0.008137398026086312
$ws.Range("E50").Value = -0.009068649678062979
$ws.Range("D51").Value = 0.00802675440491
$ws.Range("E51").Value = -0.005191350609232925
$ws.Range("D52").Value = 0.007860075619429632
$ws.Range("E52").Value = -0.01117245005257617
$ws.Range("D53").Value = 0.007217171732576779
$ws.Range("E53").Value = -0.01267893660531694
$ws.Range("D54").Value = 0.007372240071370794
$ws.Range("E54").Value = 0.006272855884472683
$ws.Range("D55").Value = 0.006733813094045308
$ws.Range("E55").Value = 0.0006575342465753309
$ws.Range("D56").Value = 0.006637239759954999
$ws.Range("E56").Value = 0.008005218216318832
$ws.Range("D57").Value = 0.006864725800114748
$ws.Range("E57").Value = -0.007911936704506517
$ws.Range("D58").Value = 0.006460918400107997
$ws.Range("E58").Value = -0.002878289473684181
$ws.Range("D59").Value = 0.005683854817858166
$ws.Range("E59").Value = -0.009330667428353756
$ws.Range("D60").Value = 0.006585484717946337
$ws.Range("E60").Value = -0.001344688480501843
$ws.Range("D61").Value = 0.005485936059302227
$ws.Range("E61").Value = -0.01246524975338548
$ws.Range("D62").Value = 0.005854813725633928
$ws.Range("E62").Value = 0.005646679214842631
$ws.Range("D63").Value = 0.005443528410812239
$ws.Range("E63").Value = 0.00535030005061099
$ws.Range("D64").Value = 0.004983439862461469
$ws.Range("E64").Value = -0.002843152740483301
$ws.Range("D65").Value = 0.004825026140723928
$ws.Range("E65").Value = -0.008320078306619294
$ws.Range("D66").Value = 0.004433321155483462
$ws.Range("E66").Value = 0.001731140555296662
$ws.Range("D67").Value = 0.004483797080484306
$ws.Range("E67").Value = -0.02172481895984224
$ws.Range("D68").Value = 0.003491103888801046
$ws.Range("E68").Value = 0.02663397311236992
$ws.Range("D69").Value = 0.00421990540240972
$ws.Range("E69").Value = -0.00657526580861767
$ws.Range("D70").Value = 0.003744861024429069
$ws.Range("E70").Value = -0.02601156069364152
$ws.Range("D71").Value = 0.00309939890356058
$ws.Range("E71").Value = 0.008904761904761971
$ws.Range("D72").Value = 0.002532996053136883
$ws.Range("E72").Value = 0.03701905335327349
$ws.Range("D73").Value = 0.002525714925553914
$ws.Range("E73").Value = 0.01610861138705455
$ws.Range("D74").Value = 0.002289914624843345
$ws.Range("E74").Value = 0.02129076612383463
$ws.Range("D75").Value = 0.001868494767575093
$ws.Range("E75").Value = 0.006529752501316377
$ws.Range("D76").Value = 0.001835532906219766
$ws.Range("E76").Value = -0.004127579737335685
$ws.Range("E77").Value = 0.003292643965357334

$ws.Protect()
